$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "20 (0.9524)"
$ws.Range("D3").Value = "1 (0.0476)"
$ws.Range("D4").Value = "4 (0.2857)"
$ws.Range("C4").Value = "10 (0.7143)"

$ws.Range("H9").Select()
